$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header "Save" in H1, copying the header formatting/style from G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add data values in H2:H4
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
